$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for columns D, E, G across all data rows to avoid Excel
# auto-converting numeric-looking / percentage-looking strings into numbers.
$ws.Range("D2:E51").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

$ws.Range("D2").Value = "262.02"; $ws.Range("E2").Value = "1.56%"; $ws.Range("G2").Value = "2"
$ws.Range("D3").Value = "27.24"; $ws.Range("E3").Value = "1.35%"; $ws.Range("G3").Value = "2"
$ws.Range("E4").Value = "2.05%"; $ws.Range("G4").Value = "2"
$ws.Range("D5").Value = "0.06069"; $ws.Range("E5").Value = "2.86%"; $ws.Range("G5").Value = "2"
$ws.Range("D6").Value = "6.713"; $ws.Range("E6").Value = "1.11%"; $ws.Range("G6").Value = "2"
$ws.Range("D7").Value = "0.8616"; $ws.Range("E7").Value = "0.99%"; $ws.Range("G7").Value = "2"
$ws.Range("D8").Value = "0.9187"; $ws.Range("E8").Value = "-2.44%"; $ws.Range("G8").Value = "2"
$ws.Range("D9").Value = "0.1411"; $ws.Range("E9").Value = "0.39%"; $ws.Range("G9").Value = "2"
$ws.Range("D10").Value = "0.04971"; $ws.Range("E10").Value = "7.16%"; $ws.Range("G10").Value = "2"
$ws.Range("D11").Value = "0.07161"; $ws.Range("E11").Value = "1.03%"; $ws.Range("G11").Value = "2"
$ws.Range("D12").Value = "0.03076"; $ws.Range("E12").Value = "-1.71%"; $ws.Range("G12").Value = "2"
$ws.Range("D13").Value = "0.09117"; $ws.Range("E13").Value = "-0.33%"; $ws.Range("G13").Value = "2"
$ws.Range("D14").Value = "0.001526"; $ws.Range("E14").Value = "0.07%"; $ws.Range("G14").Value = "2"
$ws.Range("D15").Value = "0.0006059"; $ws.Range("E15").Value = "-0.17%"; $ws.Range("G15").Value = "2"
$ws.Range("D16").Value = "0.006194"; $ws.Range("E16").Value = "2.42%"; $ws.Range("G16").Value = "2"
$ws.Range("D17").Value = "3.486"; $ws.Range("E17").Value = "-1.22%"; $ws.Range("G17").Value = "2"
$ws.Range("D18").Value = "3.166"; $ws.Range("E18").Value = "-0.71%"; $ws.Range("G18").Value = "2"
$ws.Range("D19").Value = "2.176"; $ws.Range("E19").Value = "-1.25%"; $ws.Range("G19").Value = "2"
$ws.Range("E20").Value = "2.39%"; $ws.Range("G20").Value = "2"
$ws.Range("D21").Value = "0.1289"; $ws.Range("E21").Value = "-0.85%"; $ws.Range("G21").Value = "2"
$ws.Range("D22").Value = "4.097"; $ws.Range("E22").Value = "7.34%"; $ws.Range("G22").Value = "2"
$ws.Range("D23").Value = "0.04270"; $ws.Range("E23").Value = "0.01%"; $ws.Range("G23").Value = "2"
$ws.Range("D24").Value = "0.001219"; $ws.Range("E24").Value = "0.04%"; $ws.Range("G24").Value = "2"
$ws.Range("D25").Value = "0.003931"; $ws.Range("E25").Value = "-8.44%"; $ws.Range("G25").Value = "2"
$ws.Range("D26").Value = "0.0001201"; $ws.Range("E26").Value = "0.06%"; $ws.Range("G26").Value = "2"
$ws.Range("E27").Value = "-19.26%"; $ws.Range("G27").Value = "2"
$ws.Range("G28").Value = "2"
$ws.Range("G29").Value = "2"
$ws.Range("G30").Value = "2"
$ws.Range("G31").Value = "2"
$ws.Range("G32").Value = "2"
$ws.Range("G33").Value = "2"
$ws.Range("G34").Value = "2"
$ws.Range("G35").Value = "2"
$ws.Range("G36").Value = "2"
$ws.Range("G37").Value = "2"
$ws.Range("G38").Value = "2"
$ws.Range("G39").Value = "2"
$ws.Range("D40").Value = "0.03882"; $ws.Range("E40").Value = "1.58%"; $ws.Range("G40").Value = "2"
$ws.Range("E41").Value = "1.49%"; $ws.Range("G41").Value = "2"
$ws.Range("E42").Value = "-33.26%"; $ws.Range("G42").Value = "2"
$ws.Range("D43").Value = "0.01493"; $ws.Range("E43").Value = "25.14%"; $ws.Range("G43").Value = "2"
$ws.Range("D44").Value = "0.002205"; $ws.Range("E44").Value = "16.06%"; $ws.Range("G44").Value = "2"
$ws.Range("D45").Value = "0.00005296"; $ws.Range("E45").Value = "-3.08%"; $ws.Range("G45").Value = "2"
$ws.Range("E46").Value = "0.05%"; $ws.Range("G46").Value = "2"
$ws.Range("E47").Value = "7.00%"; $ws.Range("G47").Value = "2"
$ws.Range("D48").Value = "0.1321"; $ws.Range("E48").Value = "-46.75%"; $ws.Range("G48").Value = "2"
$ws.Range("E49").Value = "0.05%"; $ws.Range("G49").Value = "2"
$ws.Range("E50").Value = "0.05%"; $ws.Range("G50").Value = "2"
$ws.Range("G51").Value = "2"

Write-Output "Updated cryptos sheet values for rows 2-51"
